$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. "Experiment results" slide (slide 11): two small text tweaks.
# ---------------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$body11 = $s11.Shapes.Item(2).TextFrame.TextRange

# 1a. Touch the trailing "/" inside the red directory run so it gets split
#     into its own run (purely cosmetic — the visible text is unchanged).
$para3 = $body11.Paragraphs(3)
$para3Text = $para3.Text
$marker = "IL.18R1/INTERVAL/r2-0.8-pruned"
$idx = $para3Text.IndexOf($marker)
if ($idx -ge 0) {
    $slashPos = $idx + $marker.Length + 1
    $slashChar = $para3.Characters($slashPos, 1)
    $slashChar.Text = $slashChar.Text
}

# 1b. Extend the "... already worked." sentence with the JAM caveat.
$para4 = $body11.Paragraphs(4)
$para4.Text = "Pipeline works on all three software but not set as a priority since both GCTA and finemap already worked. Unfortunately JAM does not select the sentinel."

# ---------------------------------------------------------------------------
# 2. Insert a new "Side results" slide right before "Other aspects".
#    Duplicating "Other aspects" keeps its layout/formatting; the original
#    copy (still slide 12) becomes "Side results" while the duplicate
#    (pushed to slide 13) keeps the old "Other aspects" content untouched.
# ---------------------------------------------------------------------------
$other = $p.Slides.Item(12)
$other.Duplicate() | Out-Null

$sideResults = $p.Slides.Item(12)
$sideResults.Shapes.Item(1).TextFrame.TextRange.Text = "Side results"

$sideBody = $sideResults.Shapes.Item(2).TextFrame.TextRange
$sideBody.Text = "sentinels() for sentinel selection is part of R/gap.`r" + `
    "cs() for credible sets is part of R/gap.`r" + `
    "With script using unpruned reference ready for GCTA and finemap, the pruned version is also OK to include JAM. The alignment of effect allele for finemap and use of R rbgen/plink2R for JAM reference file as with missing data imputation are eventually brought to control.`r" + `
    "The handling of alleles should facilitate other downstream analysis such as MR/colocalisation otherwise seen to be difficult."
